$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.288358211517334
$ws.Range("B1").Value = 2.393049955368042
$ws.Range("C1").Value = 3.177851676940918
$ws.Range("D1").Value = 3.453905820846558
$ws.Range("E1").Value = 1.062260150909424
